$wb = $excel.ActiveWorkbook

# --- Electric sheet: rename CWS_No_1..CWS_No_14 -> CWNS_No_1..CWNS_No_14 (header row, B1:O1) ---
$wsElectric = $wb.Worksheets.Item("Electric")
for ($i = 1; $i -le 14; $i++) {
    $col = $i + 1  # B=2 .. O=15
    $wsElectric.Cells.Item(1, $col).Value = "CWNS_No_$i"
}

# --- Gas sheet: rename CWS_No_1..CWS_No_7 -> CWNS_No_1..CWNS_No_7 (header row, B1:H1) ---
# and clear the stray CWS_No_8 value that used to sit in I1
$wsGas = $wb.Worksheets.Item("Gas")
for ($i = 1; $i -le 7; $i++) {
    $col = $i + 1  # B=2 .. H=8
    $wsGas.Cells.Item(1, $col).Value = "CWNS_No_$i"
}
$wsGas.Cells.Item(1, 9).Value = $null
